$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# New file identifiers / timestamps introduced by this handback run.
# ---------------------------------------------------------------------------
$oldFile = "f6bf5b5f-04f3-493a-8e89-5ae441e222d8"
$newFile1 = "2545a184-5ca5-4e30-949d-554e46b672f9"     # replaces $oldFile (row 2)
$newFile2 = "9bd0504e-e9ab-4178-8393-021487888d92"     # brand-new file (row 3)

$zhHash = "04d56327738f6a8aa55d368838971d14f0dce0a4"
$deHash = "04d56327738f6a8aa55d368838971d14f0dce0a4"
$zhHash2 = "a4e5b1c10e2638e108b9babcc9f18a0b426c3d62"
$deHash2 = "a4e5b1c10e2638e108b9babcc9f18a0b426c3d62"

# =====================================================================
# Sheet "Overview"
# =====================================================================
$ws = $wb.Worksheets.Item("Overview")
$lo = $ws.ListObjects.Item(1)

# -- row 2: rename the handed-back file + refresh the "Latest HO Xliff
#    Generate Date" column
$ws.Range("A2").Value = "$newFile1.md"
$ws.Range("B2").Value = "e2e\$newFile1.md"
$ws.Range("G2").Value = "2016-09-01 17:11:45"

# -- row 3: brand-new file handed back in this run
$lo.ListRows.Add() | Out-Null
$ws.Range("A3").Value = "$newFile2.md"
$ws.Range("C3").Value = ".md"
$ws.Range("E3").Value = "Handed back: in sync with en-US"
$ws.Range("F3").Value = "Handed back: in sync with en-US"
$ws.Range("G3").Value = "2016-09-01 17:11:45"

$ws.Hyperlinks.Add($ws.Range("B3"), `
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e1ce5889453752c0a06f869bd208d940b8a3401c/e2e/$newFile2.md", `
    "", "", "e2e\$newFile2.md") | Out-Null

# =====================================================================
# Sheet "zh-cn"
# =====================================================================
$ws = $wb.Worksheets.Item("zh-cn")
$lo = $ws.ListObjects.Item(1)

# -- row 2: rename the handed-back file + refresh timestamps
$ws.Range("A2").Value = "$newFile1.md"
$ws.Range("G2").Value = "$newFile1.$zhHash.zh-cn.xlf"
$ws.Range("H2").Value = "2016-09-01 17:11:40"
$ws.Range("I2").Value = "$newFile1.md"
$ws.Range("J2").Value = "$newFile1.$zhHash.zh-cn.xlf"
$ws.Range("K2").Value = "2016-09-01 17:11:57"

# -- row 3: brand-new file handed back in this run
$lo.ListRows.Add() | Out-Null
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Handed back: in sync with en-US"
$ws.Range("D3").Value = "e2e"
$ws.Range("E3").Value = "ht"
$ws.Range("F3").Value = "'True"
$ws.Range("G3").Value = "$newFile2.$zhHash2.zh-cn.xlf"
$ws.Range("H3").Value = "2016-09-01 17:11:40"
$ws.Range("J3").Value = "$newFile2.$zhHash2.zh-cn.xlf"
$ws.Range("K3").Value = "2016-09-01 17:11:57"
$ws.Range("L3").Value = "'"
$ws.Range("M3").Value = "'True"
$ws.Range("N3").Value = "'"
$ws.Range("O3").Value = "'False"
$ws.Range("P3").Value = "'"

$ws.Hyperlinks.Add($ws.Range("A3"), `
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/9663014cb2e850a3e027ded85c9e21b66e01a754/e2e/$newFile2.md", `
    "", "", "$newFile2.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I3"), `
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/9663014cb2e850a3e027ded85c9e21b66e01a754/e2e/$newFile2.md", `
    "", "", "$newFile2.md") | Out-Null

# =====================================================================
# Sheet "de-de"
# =====================================================================
$ws = $wb.Worksheets.Item("de-de")
$lo = $ws.ListObjects.Item(1)

# -- row 2: rename the handed-back file + refresh timestamps
$ws.Range("A2").Value = "$newFile1.md"
$ws.Range("G2").Value = "$newFile1.$deHash.de-de.xlf"
$ws.Range("H2").Value = "2016-09-01 17:11:45"
$ws.Range("I2").Value = "$newFile1.md"
$ws.Range("J2").Value = "$newFile1.$deHash.de-de.xlf"
$ws.Range("K2").Value = "2016-09-01 17:12:15"

# -- row 3: brand-new file handed back in this run
$lo.ListRows.Add() | Out-Null
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Handed back: in sync with en-US"
$ws.Range("D3").Value = "e2e"
$ws.Range("E3").Value = "ht"
$ws.Range("F3").Value = "'True"
$ws.Range("G3").Value = "$newFile2.$deHash2.de-de.xlf"
$ws.Range("H3").Value = "2016-09-01 17:11:45"
$ws.Range("J3").Value = "$newFile2.$deHash2.de-de.xlf"
$ws.Range("K3").Value = "2016-09-01 17:12:15"
$ws.Range("L3").Value = "'"
$ws.Range("M3").Value = "'True"
$ws.Range("N3").Value = "'"
$ws.Range("O3").Value = "'False"
$ws.Range("P3").Value = "'"

$ws.Hyperlinks.Add($ws.Range("A3"), `
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/6b8a101cd01f43fa8be0ff135809b9ceeba209fc/e2e/$newFile2.md", `
    "", "", "$newFile2.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I3"), `
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/6b8a101cd01f43fa8be0ff135809b9ceeba209fc/e2e/$newFile2.md", `
    "", "", "$newFile2.md") | Out-Null

Write-Output "Report generated for handback"
